# Update "想去人数" (interest count) figures that were refreshed by the
# upstream scraper run, for both the "展览" sheet and the "全部类型" sheet
# (which mirrors the same events).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1580
$ws1.Range("F5").Value  = 173
$ws1.Range("F7").Value  = 5211
$ws1.Range("F9").Value  = 771
$ws1.Range("F10").Value = 1059
$ws1.Range("F15").Value = 30
$ws1.Range("F16").Value = 6571
$ws1.Range("F19").Value = 142
$ws1.Range("F20").Value = 170
$ws1.Range("F22").Value = 1017
$ws1.Range("F23").Value = 15712
$ws1.Range("F28").Value = 110
$ws1.Range("F29").Value = 11153
$ws1.Range("F30").Value = 795
$ws1.Range("F32").Value = 270

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1580
$ws4.Range("F5").Value  = 173
$ws4.Range("F8").Value  = 5211
$ws4.Range("F10").Value = 771
$ws4.Range("F12").Value = 1059
$ws4.Range("F18").Value = 30
$ws4.Range("F19").Value = 6571
$ws4.Range("F22").Value = 142
$ws4.Range("F23").Value = 170
$ws4.Range("F26").Value = 1017
$ws4.Range("F27").Value = 15712
$ws4.Range("F32").Value = 110
$ws4.Range("F34").Value = 11153
$ws4.Range("F35").Value = 795
$ws4.Range("F37").Value = 270
